# Updated cryptos list with refreshed prices / volume(1h) figures.
# Source data column D stores prices as text (it includes thousands
# separators like "96.346.97" that aren't valid numbers), so force the
# whole data range to Text format first - this preserves values such as
# "242.00" or "1.00" exactly as strings instead of Excel normalising them
# to numeric 242 / 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '96.346.97'
$ws.Range("D3").Value = '3.659.32'
$ws.Range("D5").Value = '242.00'
$ws.Range("D7").Value = '661.01'
$ws.Range("D11").Value = '3.657.58'
$ws.Range("D12").Value = '44.81'
$ws.Range("D14").Value = '6.67'
$ws.Range("D15").Value = '4.338.25'
$ws.Range("D17").Value = '96.188.15'
$ws.Range("D18").Value = '8.92'
$ws.Range("D19").Value = '3.652.83'
$ws.Range("D20").Value = '12.74'
$ws.Range("D21").Value = '18.26'
$ws.Range("D23").Value = '520.25'
$ws.Range("D24").Value = '3.43'
$ws.Range("D27").Value = '102.15'
$ws.Range("D28").Value = '12.97'
$ws.Range("D31").Value = '3.04'
$ws.Range("D32").Value = '1.00'
$ws.Range("D35").Value = '32.95'
$ws.Range("D38").Value = '627.66'
$ws.Range("D39").Value = '45.22'
$ws.Range("D40").Value = '8.72'
$ws.Range("D46").Value = '0.0457'
$ws.Range("D47").Value = '0.444'
$ws.Range("D49").Value = '23.60'
$ws.Range("D50").Value = '8.55'
$ws.Range("D51").Value = '3.57'

$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("E3").Value = '  +1.77%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("E6").Value = '  +14.70%  '
$ws.Range("E7").Value = '  +0.77%  '
$ws.Range("E8").Value = '  +2.35%  '
$ws.Range("E9").Value = '  +2.52%  '
$ws.Range("E10").Value = '  +0.00%  '
$ws.Range("E11").Value = '  +1.81%  '
$ws.Range("E12").Value = '  +2.36%  '
$ws.Range("E13").Value = '  +0.24%  '
$ws.Range("E14").Value = '  +3.35%  '
$ws.Range("E15").Value = '  +1.78%  '
$ws.Range("E16").Value = '  +5.27%  '
$ws.Range("E17").Value = '  -0.93%  '
$ws.Range("E18").Value = '  +15.08%  '
$ws.Range("E19").Value = '  +1.87%  '
$ws.Range("E21").Value = '  +0.71%  '
$ws.Range("E22").Value = '  +0.80%  '
$ws.Range("E23").Value = '  +1.76%  '
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("E26").Value = '  -0.64%  '
$ws.Range("E27").Value = '  +4.97%  '
$ws.Range("E28").Value = '  -1.69%  '
$ws.Range("E29").Value = '  +9.52%  '
$ws.Range("E30").Value = '  +6.02%  '
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("E35").Value = '  +4.09%  '
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("E37").Value = '  +2.57%  '
$ws.Range("E38").Value = '  +0.58%  '
$ws.Range("E39").Value = '  +35.57%  '
$ws.Range("E40").Value = '  -1.40%  '
$ws.Range("E41").Value = '  +4.71%  '
$ws.Range("E42").Value = '  +4.35%  '
$ws.Range("E43").Value = '  +5.06%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("E45").Value = '  +7.69%  '
$ws.Range("E46").Value = '  +5.33%  '
$ws.Range("E47").Value = '  +23.48%  '
$ws.Range("E48").Value = '  -1.36%  '
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("E50").Value = '  +2.61%  '
$ws.Range("E51").Value = '  +1.14%  '
